$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# --- New paragraph 1: numbered list item continuing the "hk"/numId=2 list ---
$p1xml = "<w:p $wns>" +
    "<w:pPr>" +
        "<w:pStyle w:val='hk'/>" +
        "<w:numPr><w:ilvl w:val='0'/><w:numId w:val='2'/></w:numPr>" +
        "<w:shd w:val='clear' w:color='auto' w:fill='FFFFFF'/>" +
        "<w:spacing w:before='480' w:beforeAutospacing='0' w:after='0' w:afterAutospacing='0'/>" +
        "<w:jc w:val='both'/>" +
        "<w:rPr>" +
            "<w:rFonts w:asciiTheme='minorHAnsi' w:hAnsiTheme='minorHAnsi' w:cstheme='minorHAnsi'/>" +
            "<w:spacing w:val='-1'/>" +
            "<w:szCs w:val='32'/>" +
        "</w:rPr>" +
    "</w:pPr>" +
    "<w:r>" +
        "<w:rPr>" +
            "<w:rFonts w:asciiTheme='minorHAnsi' w:hAnsiTheme='minorHAnsi' w:cstheme='minorHAnsi'/>" +
            "<w:spacing w:val='-1'/>" +
            "<w:sz w:val='28'/>" +
            "<w:szCs w:val='32'/>" +
        "</w:rPr>" +
        "<w:t>In case of error the whole application crashed due to high dependency of services over each other.</w:t>" +
    "</w:r>" +
"</w:p>"

# --- New paragraph 2: empty "hk" paragraph (indented, not numbered) holding the _GoBack bookmark ---
$p2xml = "<w:p $wns>" +
    "<w:pPr>" +
        "<w:pStyle w:val='hk'/>" +
        "<w:shd w:val='clear' w:color='auto' w:fill='FFFFFF'/>" +
        "<w:spacing w:before='480' w:beforeAutospacing='0' w:after='0' w:afterAutospacing='0'/>" +
        "<w:ind w:left='360'/>" +
        "<w:jc w:val='both'/>" +
        "<w:rPr>" +
            "<w:rFonts w:asciiTheme='minorHAnsi' w:hAnsiTheme='minorHAnsi' w:cstheme='minorHAnsi'/>" +
            "<w:spacing w:val='-1'/>" +
            "<w:szCs w:val='32'/>" +
        "</w:rPr>" +
    "</w:pPr>" +
    "<w:bookmarkStart w:id='0' w:name='_GoBack'/>" +
    "<w:bookmarkEnd w:id='0'/>" +
"</w:p>"

# Find the paragraph that ends the "Drawbacks of Monolithic Architecture" list
# ("All the services and parts of monolithic architecture ... separately.").
$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $pp = $d.Paragraphs($i)
    if ($pp.Range.Text -like "*All the services and parts of monolithic architecture*") {
        $anchor = $pp
    }
}

# Drop the pre-existing "_GoBack" bookmark up front so re-adding it below (inside the
# new second paragraph) doesn't leave two bookmarks of the same name behind.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# Insert a placeholder paragraph right after the anchor, then stamp it with the
# full target OOXML for the first new paragraph.
$anchor.Range.InsertParagraphAfter() | Out-Null
$placeholder1 = $anchor.Next()
$placeholder1.Range.InsertXML($p1xml) | Out-Null

# Re-resolve the paragraph we just wrote (InsertXML reseats ranges) and insert the
# second placeholder paragraph after it, then stamp it with the second paragraph's XML.
$newPara1 = $anchor.Next()
$newPara1.Range.InsertParagraphAfter() | Out-Null
$placeholder2 = $newPara1.Next()
$placeholder2.Range.InsertXML($p2xml) | Out-Null
